# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Vega Modelo de Temuco" / Espinaca
# as row 215, shifting the existing rows 215-227 down to 216-228.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 215..227 down by one to make room for the new record.
$ws.Rows("215:215").Insert()

# Populate the newly inserted row 215 with the new weekly record.
$ws.Cells.Item(215, 1).Value = 10
$ws.Cells.Item(215, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(215, 3).Value = "La Araucanía"
$ws.Cells.Item(215, 4).Value = 44931
$ws.Cells.Item(215, 5).Value = 9
$ws.Cells.Item(215, 6).Value = 100112012
$ws.Cells.Item(215, 7).Value = "Espinaca"
$ws.Cells.Item(215, 8).Value = "Sin especificar"
$ws.Cells.Item(215, 9).Value = "Primera"
$ws.Cells.Item(215, 10).Value = 90
$ws.Cells.Item(215, 11).Value = 10000
$ws.Cells.Item(215, 12).Value = 10000
$ws.Cells.Item(215, 13).Value = 10000
$ws.Cells.Item(215, 14).Value = "`$/docena de atados"
$ws.Cells.Item(215, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(215, 16).Value = 3333
$ws.Cells.Item(215, 17).Value = 3
$ws.Cells.Item(215, 18).Value = "Hortaliza"
